# Adds address info to the batch/lot notification letter template:
#   - switches every paragraph to single line spacing
#     (<w:spacing w:line="240" w:lineRule="auto"/>)
#   - replaces the hard-coded date with a {{formatted_date}} placeholder
#   - replaces the hard-coded bank/trust owner line with a {{owner}} placeholder
#   - underlines the original-batch-address placeholders / section titles
#   - indents the first line of the intro paragraph
#   - removes the decorative letterhead image that used to live in the header
#     (now that {{owner}}/address lines make the header redundant)

$d = $word.ActiveDocument

# wdLineSpaceSingle = 0  -> serializes as <w:spacing w:line="240" w:lineRule="auto"/>
$wdLineSpaceSingle = 0
# wdUnderlineSingle = 1
$wdUnderlineSingle = 1

# --- 1. Give every paragraph single line spacing -------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $d.Paragraphs.Item($i).Format.LineSpacingRule = $wdLineSpaceSingle
}

# --- 2. Replace the hard-coded date with the {{formatted_date}} token ----
$d.Paragraphs.Item(6).Range.Text = "PLAYA DEL CARMEN, Q.ROO, A {{formatted_date}}"

# --- 3. Replace the hard-coded owner/trust line with {{owner}} -----------
$d.Paragraphs.Item(8).Range.Text = "{{owner}}"

# --- 4. Underline the original batch address placeholder under the owner -
$d.Paragraphs.Item(9).Range.Font.Underline = $wdUnderlineSingle

# --- 5. Indent the first line of the "En respuesta..." paragraph --------
$d.Paragraphs.Item(11).Format.FirstLineIndent = 36

# --- 6. Underline "{{original_batch_address}}: " under LOTE ORIGINAL ----
$d.Paragraphs.Item(14).Range.Font.Underline = $wdUnderlineSingle

# --- 7. Underline the "LOTES RESULTANTES:" heading -----------------------
$d.Paragraphs.Item(24).Range.Font.Underline = $wdUnderlineSingle

# --- 8. Remove the decorative letterhead image from the page header -----
$header = $d.Sections.Item(1).Headers.Item(1)
if ($header.Shapes.Count -gt 0) {
    for ($i = $header.Shapes.Count; $i -ge 1; $i--) {
        $header.Shapes.Item($i).Delete()
    }
}
